$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 3 new rows of master data (regcntr_id 10005, usr_id 110033-110035)
$newRows = @(
    @(10005, 110033, 10005),
    @(10005, 110034, 10005),
    @(10005, 110035, 10005)
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
}

# Select the rows below the newly added data, mirroring the post-edit UI state
$ws.Rows("37:1048576").Select() | Out-Null
